$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3234.125
$ws.Range("I40").Value = 2328
$ws.Range("J40").Value = 3777.8
$ws.Range("K40").Value = 2328
$ws.Range("L40").Value = 3777.8
$ws.Range("M40").Value = -2153
$ws.Range("N40").Value = -4127.8
$ws.Range("H76").Value = 20331.834
$ws.Range("I76").Value = 19726.637
$ws.Range("J76").Value = 21282.857
$ws.Range("K76").Value = 19726.637
$ws.Range("L76").Value = 21282.857
$ws.Range("M76").Value = -19411.637
$ws.Range("N76").Value = -21912.857
$ws.Range("H79").Value = 20331.834
$ws.Range("I79").Value = 19726.637
$ws.Range("J79").Value = 21282.857
$ws.Range("K79").Value = 19726.637
$ws.Range("L79").Value = 21282.857
$ws.Range("M79").Value = -18634.637
$ws.Range("N79").Value = -23466.857
$ws.Range("H80").Value = 788
$ws.Range("I80").Value = 425
$ws.Range("J80").Value = 854
$ws.Range("K80").Value = 1275
$ws.Range("L80").Value = 2562
$ws.Range("M80").Value = -277
$ws.Range("N80").Value = -4558
$ws.Range("H83").Value = 788
$ws.Range("I83").Value = 425
$ws.Range("J83").Value = 854
$ws.Range("K83").Value = 3825
$ws.Range("L83").Value = 7686
$ws.Range("M83").Value = 1167
$ws.Range("N83").Value = -17670
$ws.Range("H138").Value = 3484.1042
$ws.Range("I138").Value = 1460.8334
$ws.Range("J138").Value = 4698.067
$ws.Range("K138").Value = 4382.5002
$ws.Range("L138").Value = 14094.201
$ws.Range("M138").Value = 757.4997999999996
$ws.Range("N138").Value = -24374.201

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5209.909
$ws.Range("I32").Value = 1465.8649
$ws.Range("K32").Value = 1465.8649
$ws.Range("M32").Value = -1178.8649
$ws.Range("H62").Value = 50249
$ws.Range("J62").Value = 50249
$ws.Range("L62").Value = 50249
$ws.Range("N62").Value = -51497
$ws.Range("H63").Value = 2423.1428
$ws.Range("I63").Value = 2695.2
$ws.Range("J63").Value = 1743
$ws.Range("K63").Value = 2695.2
$ws.Range("L63").Value = 1743
$ws.Range("M63").Value = -2009.2
$ws.Range("N63").Value = -3115
$ws.Range("H65").Value = 50249
$ws.Range("J65").Value = 50249
$ws.Range("L65").Value = 150747
$ws.Range("N65").Value = -156987
$ws.Range("H66").Value = 2423.1428
$ws.Range("I66").Value = 2695.2
$ws.Range("J66").Value = 1743
$ws.Range("K66").Value = 13476
$ws.Range("L66").Value = 8715
$ws.Range("M66").Value = -10044
$ws.Range("N66").Value = -15579
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H88").Value = 10286.538
$ws.Range("I88").Value = 17638.834
$ws.Range("J88").Value = 3984.5715
$ws.Range("K88").Value = 17638.834
$ws.Range("L88").Value = 3984.5715
$ws.Range("M88").Value = -17232.834
$ws.Range("N88").Value = -4796.5715
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H91").Value = 10286.538
$ws.Range("I91").Value = 17638.834
$ws.Range("J91").Value = 3984.5715
$ws.Range("K91").Value = 17638.834
$ws.Range("L91").Value = 3984.5715
$ws.Range("M91").Value = -16234.834
$ws.Range("N91").Value = -6792.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5666.3335
$ws.Range("J20").Value = 5999.5
$ws.Range("L20").Value = 5999.5
$ws.Range("N20").Value = -6493.5
$ws.Range("H76").Value = 35535.332
$ws.Range("J76").Value = 37642.6
$ws.Range("L76").Value = 37642.6
$ws.Range("N76").Value = -38272.6
$ws.Range("H79").Value = 35535.332
$ws.Range("J79").Value = 37642.6
$ws.Range("L79").Value = 37642.6
$ws.Range("N79").Value = -39826.6
$ws.Range("H107").Value = 1635.5454
$ws.Range("I107").Value = 1635.5454
$ws.Range("K107").Value = 1635.5454
$ws.Range("M107").Value = 284.4546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2037.0869
$ws.Range("I58").Value = 2337.2727
$ws.Range("J58").Value = 1761.9166
$ws.Range("K58").Value = 2337.2727
$ws.Range("L58").Value = 1761.9166
$ws.Range("M58").Value = -2134.2727
$ws.Range("N58").Value = -2167.9166
$ws.Range("H62").Value = 6097.0835
$ws.Range("I62").Value = 5962.5
$ws.Range("K62").Value = 5962.5
$ws.Range("M62").Value = -5338.5
$ws.Range("H65").Value = 6097.0835
$ws.Range("I65").Value = 5962.5
$ws.Range("K65").Value = 29812.5
$ws.Range("M65").Value = -26692.5
$ws.Range("H132").Value = 2418.5454
$ws.Range("I132").Value = 2391.25
$ws.Range("K132").Value = 7173.75
$ws.Range("M132").Value = -4643.75
$ws.Range("H136").Value = 2037.0869
$ws.Range("I136").Value = 2337.2727
$ws.Range("J136").Value = 1761.9166
$ws.Range("K136").Value = 7011.8181
$ws.Range("L136").Value = 5285.7498
$ws.Range("M136").Value = -4461.8181
$ws.Range("N136").Value = -10385.7498

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 47514
$ws.Range("J15").Value = 47514
$ws.Range("L15").Value = 47514
$ws.Range("N15").Value = -48090
$ws.Range("H70").Value = 7627.4
$ws.Range("I70").Value = 7571.875
$ws.Range("K70").Value = 7571.875
$ws.Range("M70").Value = -7301.875
$ws.Range("H73").Value = 7627.4
$ws.Range("I73").Value = 7571.875
$ws.Range("K73").Value = 7571.875
$ws.Range("M73").Value = -6635.875
$ws.Range("H80").Value = 6665.4546
$ws.Range("I80").Value = 3837
$ws.Range("J80").Value = 8281.714
$ws.Range("K80").Value = 3837
$ws.Range("L80").Value = 8281.714
$ws.Range("M80").Value = -2839
$ws.Range("N80").Value = -10277.714
$ws.Range("H81").Value = 47514
$ws.Range("J81").Value = 47514
$ws.Range("L81").Value = 47514
$ws.Range("N81").Value = -49510
$ws.Range("H83").Value = 6665.4546
$ws.Range("I83").Value = 3837
$ws.Range("J83").Value = 8281.714
$ws.Range("K83").Value = 19185
$ws.Range("L83").Value = 41408.57
$ws.Range("M83").Value = -14193
$ws.Range("N83").Value = -51392.57
$ws.Range("H84").Value = 47514
$ws.Range("J84").Value = 47514
$ws.Range("L84").Value = 142542
$ws.Range("N84").Value = -152526
$ws.Range("H87").Value = 69999
$ws.Range("J87").Value = 69999
$ws.Range("L87").Value = 69999
$ws.Range("N87").Value = -72495
$ws.Range("H90").Value = 69999
$ws.Range("J90").Value = 69999
$ws.Range("L90").Value = 209997
$ws.Range("N90").Value = -222477

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2793.077
$ws.Range("I22").Value = 2881
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 2881
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -2586
$ws.Range("N22").Value = -3090
$ws.Range("H27").Value = 2793.077
$ws.Range("I27").Value = 2881
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 2881
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -2774
$ws.Range("N27").Value = -2714
$ws.Range("H68").Value = 3175.7693
$ws.Range("I68").Value = 2844.818
$ws.Range("J68").Value = 4996
$ws.Range("K68").Value = 2844.818
$ws.Range("L68").Value = 4996
$ws.Range("M68").Value = -2095.818
$ws.Range("N68").Value = -6494
$ws.Range("H71").Value = 3175.7693
$ws.Range("I71").Value = 2844.818
$ws.Range("J71").Value = 4996
$ws.Range("K71").Value = 14224.09
$ws.Range("L71").Value = 24980
$ws.Range("M71").Value = -10480.09
$ws.Range("N71").Value = -32468
